$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '24.435.34'
Set-TextValue 'E2' '  -1.19%  '
Set-TextValue 'D3' '1.656.66'
Set-TextValue 'E3' '  -2.78%  '
Set-TextValue 'D4' '1.006'
Set-TextValue 'E4' '  +0.29%  '
Set-TextValue 'D5' '307.68'
Set-TextValue 'E5' '  -0.35%  '
Set-TextValue 'E6' '  +0.36%  '
Set-TextValue 'D7' '0.3608'
Set-TextValue 'E7' '  -3.41%  '
Set-TextValue 'D8' '47.40'
Set-TextValue 'E8' '  -3.47%  '
Set-TextValue 'D9' '0.3250'
Set-TextValue 'E9' '  -5.50%  '
Set-TextValue 'E10' '  -6.69%  '
Set-TextValue 'D11' '0.06993'
Set-TextValue 'E11' '  -6.42%  '
Set-TextValue 'D12' '1.002'
Set-TextValue 'E12' '  +0.22%  '
Set-TextValue 'D13' '5.887'
Set-TextValue 'E13' '  -5.64%  '
Set-TextValue 'D14' '19.32'
Set-TextValue 'E14' '  -7.64%  '
Set-TextValue 'D15' '1.653.19'
Set-TextValue 'E15' '  -3.19%  '
Set-TextValue 'D16' '6.556'
Set-TextValue 'E16' '  -5.83%  '
Set-TextValue 'E17' '  -7.45%  '
Set-TextValue 'D18' '0.06554'
Set-TextValue 'E18' '  -2.45%  '
Set-TextValue 'D20' '76.56'
Set-TextValue 'E20' '  -9.32%  '
Set-TextValue 'D21' '5.929'
Set-TextValue 'E21' '  -6.28%  '
Set-TextValue 'E22' '  -8.95%  '
Set-TextValue 'D23' '12.44'
Set-TextValue 'E23' '  -4.56%  '
Set-TextValue 'D24' '24.433.95'
Set-TextValue 'E24' '  -1.21%  '
Set-TextValue 'D25' '2.464'
Set-TextValue 'E25' '  +0.95%  '
Set-TextValue 'D26' '2.303'
Set-TextValue 'E26' '  -16.53%  '
Set-TextValue 'D27' '146.97'
Set-TextValue 'E27' '  -2.06%  '
Set-TextValue 'D28' '18.47'
Set-TextValue 'E28' '  -8.68%  '
Set-TextValue 'D29' '1.837.27'
Set-TextValue 'E29' '  -3.12%  '
Set-TextValue 'D30' '123.58'
Set-TextValue 'E30' '  -5.88%  '
Set-TextValue 'D31' '1.171'
Set-TextValue 'E31' '  -0.57%  '
Set-TextValue 'D32' '3.974'
Set-TextValue 'E32' '  -4.86%  '
Set-TextValue 'D33' '5.627'
Set-TextValue 'E33' '  -16.77%  '
Set-TextValue 'D34' '1.713'
Set-TextValue 'E34' '  -4.64%  '
Set-TextValue 'D35' '0.08381'
Set-TextValue 'E35' '  -4.99%  '
Set-TextValue 'D36' '12.28'
Set-TextValue 'E36' '  -10.00%  '
Set-TextValue 'D37' '5.190'
Set-TextValue 'E37' '  -6.15%  '
Set-TextValue 'D38' '0.06023'
Set-TextValue 'E38' '  -8.52%  '
Set-TextValue 'E39' '  -8.15%  '
Set-TextValue 'B40' 'Algorand'
Set-TextValue 'C40' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D40' '0.2053'
Set-TextValue 'E40' '  -7.43%  '
Set-TextValue 'B41' 'TrustWalletToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D41' '1.201'
Set-TextValue 'E41' '  -5.68%  '
Set-TextValue 'D42' '8.183'
Set-TextValue 'E42' '  -8.87%  '
Set-TextValue 'E43' '  +0.48%  '
Set-TextValue 'D44' '0.5884'
Set-TextValue 'E44' '  -8.62%  '
Set-TextValue 'D45' '3.737'
Set-TextValue 'E45' '  -2.01%  '
Set-TextValue 'D46' '12.69'
Set-TextValue 'E46' '  -8.62%  '
Set-TextValue 'D47' '0.5585'
Set-TextValue 'E47' '  -8.58%  '
Set-TextValue 'D48' '122.00'
Set-TextValue 'E48' '  -5.68%  '
Set-TextValue 'D49' '1.929'
Set-TextValue 'E49' '  -8.87%  '
Set-TextValue 'D50' '0.06893'
Set-TextValue 'E50' '  -5.51%  '
Set-TextValue 'E51' '  -6.55%  '
